$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (column C) date for every existing data row
#    (rows 2..440) from 2023-10-06 (45205) to 2023-10-07 (45206).
for ($r = 2; $r -le 440; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}

# 2. Row 440 gains an explicit row height (matches the other data rows).
$ws.Rows.Item(440).RowHeight = 15

# 3. Append two new records (rows 441 and 442).
$ws.Range("A441").Value = "A 48057-2023"
$ws.Range("B441").Value = 45204
$ws.Range("B441").NumberFormat = "YYYY-MM-DD"
$ws.Range("C441").Value = 45206
$ws.Range("C441").NumberFormat = "YYYY-MM-DD"
$ws.Range("D441").Value = "NORRBOTTENS LÄN"
$ws.Range("E441").Value = "KALIX"
$ws.Range("G441").Value = 0.8
$ws.Range("H441").Value = 0
$ws.Range("I441").Value = 0
$ws.Range("J441").Value = 0
$ws.Range("K441").Value = 0
$ws.Range("L441").Value = 0
$ws.Range("M441").Value = 0
$ws.Range("N441").Value = 0
$ws.Range("O441").Value = 0
$ws.Range("P441").Value = 0
$ws.Range("Q441").Value = 0
$ws.Range("R441").WrapText = $true
$ws.Rows.Item(441).RowHeight = 15

$ws.Range("A442").Value = "A 48063-2023"
$ws.Range("B442").Value = 45204
$ws.Range("B442").NumberFormat = "YYYY-MM-DD"
$ws.Range("C442").Value = 45206
$ws.Range("C442").NumberFormat = "YYYY-MM-DD"
$ws.Range("D442").Value = "NORRBOTTENS LÄN"
$ws.Range("E442").Value = "KALIX"
$ws.Range("G442").Value = 2.1
$ws.Range("H442").Value = 0
$ws.Range("I442").Value = 0
$ws.Range("J442").Value = 0
$ws.Range("K442").Value = 0
$ws.Range("L442").Value = 0
$ws.Range("M442").Value = 0
$ws.Range("N442").Value = 0
$ws.Range("O442").Value = 0
$ws.Range("P442").Value = 0
$ws.Range("Q442").Value = 0
$ws.Range("R442").WrapText = $true
